$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G3 previously held "MSRP_2000_CHOICE_T" - rename to "MSRP_2000_CHOICE"
# (this text is shared with G2, so the duplicate shared-string entry goes away)
$ws.Range("G3").Value = "MSRP_2000_CHOICE"

# Update the active selection to match the new authoring session
$ws.Range("P7").Select()
